$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position (best effort; not exposed through this COM surface) ---
# xWindow="930" -> "1860" in the workbook's bookViews. There is no reachable
# Window/Application property on this runtime that maps onto that attribute,
# so this is intentionally left alone.

# --- Shared-string juggling so the final table matches the authored order ---
# "offset" (the lone occupant of its shared-string slot) is retargeted to
# "increment" first; a brand new cell (K9) then also becomes "increment" and
# picks up that same slot; finally K7 is renamed again to "idle offset",
# which forks off a fresh slot (since K9 still needs "increment") that ends
# up last in the table - matching increment@23 ... idle offset@35.
$ws.Range("K7").Value = "increment"
$ws.Range("K9").Value = "increment"

# --- New "analog waveform generator" register rows (8-11) ---
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = "BOOL"
$ws.Range("J8").Value = "W"
$ws.Range("K8").Value = "active?"

$ws.Range("H9").Value = 2
$ws.Range("I9").Value = "U32"
$ws.Range("J9").Value = "W"
# K9 already set to "increment" above

$ws.Range("H10").Value = 3
$ws.Range("I10").Value = "U16"
$ws.Range("J10").Value = "W"
$ws.Range("K10").Value = "number of elements"

$ws.Range("H11").Value = "4..1026"
$ws.Range("I11").Value = "I16"
$ws.Range("J11").Value = "W"
$ws.Range("K11").Value = "LUT"

# K7 gets its final text last, after every other new string referencing
# "increment" is already in place.
$ws.Range("K7").Value = "idle offset"

# --- Selection moves from A9 to O9 ---
$ws.Range("O9").Select()
